$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-14: matches that have now been played - fill in results (xG_home,
# xG_away, goals_home, goals_away) and fix up home/away team ordering that
# shifted as the fixture list moved forward.
$ws.Range("B10").Value = "Lazio"
$ws.Range("C10").Value = "Udinese"

$ws.Range("B11").Value = "Torino"
$ws.Range("C11").Value = "Udinese"

$ws.Range("B12").Value = "Udinese"
$ws.Range("C12").Value = "Crotone"

$ws.Range("B13").Value = "Cagliari"
$ws.Range("C13").Value = "Udinese"

$ws.Range("B14").Value = "Udinese"
$ws.Range("C14").Value = "Benevento"

# Remaining still-pending fixtures (rows 15-19) shift up one slot.
$ws.Range("B15").Value = "Juventus"
$ws.Range("C15").Value = "Udinese"

$ws.Range("B16").Value = "Bologna"
$ws.Range("C16").Value = "Udinese"

$ws.Range("B17").Value = "Udinese"
$ws.Range("C17").Value = "Napoli"

$ws.Range("B18").Value = "Sampdoria"
$ws.Range("C18").Value = "Udinese"

$ws.Range("B19").Value = "Udinese"
$ws.Range("C19").Value = "Atalanta"

# New result data (xG/goals) for rows 10-14 - stored as text in the source
# workbook, so force a text number format before writing, then drop the
# number format again so no stray cell style lingers on save.
$resultsRange = $ws.Range("D10:G14")
$resultsRange.NumberFormat = "@"

$ws.Range("D10").Value = "1.22788"
$ws.Range("E10").Value = "1.42863"
$ws.Range("F10").Value = "1"
$ws.Range("G10").Value = "3"

$ws.Range("D11").Value = "1.07631"
$ws.Range("E11").Value = "1.47406"
$ws.Range("F11").Value = "2"
$ws.Range("G11").Value = "3"

$ws.Range("D12").Value = "1.10391"
$ws.Range("E12").Value = "0.123143"
$ws.Range("F12").Value = "0"
$ws.Range("G12").Value = "0"

$ws.Range("D13").Value = "0.628497"
$ws.Range("E13").Value = "1.43641"
$ws.Range("F13").Value = "1"
$ws.Range("G13").Value = "1"

$ws.Range("D14").Value = "2.4683"
$ws.Range("E14").Value = "0.338064"
$ws.Range("F14").Value = "0"
$ws.Range("G14").Value = "2"

$resultsRange.Style = "Normal"
